$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Persentase Tingkat Pendidikan"
